# Fixed import of PI and Investment instrument uniqueness
#
# The "Folio No" column (previously column G) was removed from the import
# template entirely. Deleting the whole column shifts every column to its
# right (Instrument, Currency, Sector, Investment Domicile *, Custom Field 1)
# one position to the left, which also collapses the duplicate "Instrument"
# shared-string entries that the import had produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "Folio No" header (row 1) / blank cells (rows 2-3).
# Deleting the entire column shifts H:L left to G:K.
$ws.Range("G1").EntireColumn.Delete()
